$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.003.73"
$ws.Range("E2").Value = "  -8.04%  "
$ws.Range("D3").Value = "1.420.18"
$ws.Range("E3").Value = "  -7.76%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'273.96"
$ws.Range("E6").Value = "  -5.50%  "
$ws.Range("D7").Value = "'0.3717"
$ws.Range("E7").Value = "  -4.12%  "
$ws.Range("D8").Value = "'0.3074"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").Value = "'39.69"
$ws.Range("E9").Value = "  -7.90%  "
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "'0.06605"
$ws.Range("E11").Value = "  -8.28%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'5.421"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "'17.14"
$ws.Range("E14").Value = "  -7.67%  "
$ws.Range("D15").Value = "'6.168"
$ws.Range("E15").Value = "  -6.51%  "
$ws.Range("D16").Value = "1.421.11"
$ws.Range("E16").Value = "  -7.75%  "
$ws.Range("D17").Value = "'0.00001008"
$ws.Range("E17").Value = "  -9.15%  "
$ws.Range("D18").Value = "'0.05813"
$ws.Range("E18").Value = "  -11.80%  "
$ws.Range("D19").Value = "'74.64"
$ws.Range("E19").Value = "  -10.38%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'5.647"
$ws.Range("E21").Value = "  -7.99%  "
$ws.Range("D22").Value = "'14.50"
$ws.Range("E22").Value = "  -5.78%  "
$ws.Range("D23").Value = "'11.06"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "'2.340"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "20.019.94"
$ws.Range("E25").Value = "  -7.95%  "
$ws.Range("D26").Value = "'2.301"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'138.88"
$ws.Range("E27").Value = "  -5.41%  "
$ws.Range("D28").Value = "'16.90"
$ws.Range("E28").Value = "  -7.98%  "
$ws.Range("D29").Value = "1.581.00"
$ws.Range("E29").Value = "  -7.82%  "
$ws.Range("D30").Value = "'109.03"
$ws.Range("D31").Value = "'3.812"
$ws.Range("E31").Value = "  -21.31%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.431"
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.8889"
$ws.Range("E33").Value = "  -8.30%  "
$ws.Range("D34").Value = "'0.07733"
$ws.Range("E34").Value = "  -5.51%  "
$ws.Range("D35").Value = "'8.442"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").Value = "'11.33"
$ws.Range("E36").Value = "  +6.27%  "
$ws.Range("D37").Value = "'4.786"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05704"
$ws.Range("E38").Value = "  -6.05%  "
$ws.Range("B39").Value = "Frax"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D39").Value = "'1.000"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'0.1924"
$ws.Range("E40").Value = "  -5.51%  "
$ws.Range("E41").Value = "  -7.68%  "
$ws.Range("D42").Value = "'1.091"
$ws.Range("E42").Value = "  -8.36%  "
$ws.Range("D43").Value = "'1.268"
$ws.Range("E43").Value = "  -14.59%  "
$ws.Range("D44").Value = "'0.5331"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.30"
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.537"
$ws.Range("E46").Value = "  -5.56%  "
$ws.Range("D47").Value = "'0.5137"
$ws.Range("E47").Value = "  -6.89%  "
$ws.Range("D48").Value = "'1.802"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "'109.56"
$ws.Range("E49").Value = "  -7.02%  "
$ws.Range("E50").Value = "  -8.29%  "
$ws.Range("D51").Value = "'0.9993"
$ws.Range("E51").Value = "  -0.16%  "
